$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (new C and D) before the old column C,
# shifting the old "date" column (C) to E.
$ws.Columns("C:D").Insert()

# Copy the (now shifted) original column C formatting from E onto the
# two freshly inserted columns so they pick up the same number format /
# style as the rest of the date columns.
$ws.Range("E1:E27").Copy()
$ws.Range("C1:D27").PasteSpecial(-4122)

# Keep all three date columns at the same fixed width the original
# column C used (8 characters).
$ws.Columns("C:E").ColumnWidth = 7.083333333333333

# The old B1 header ("Jun_13") now needs to move to the newly inserted
# D column before we overwrite B1 with the newest date.
$ws.Range("D1").Value = "Jun_13"

# New header dates for the two newly-inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill in the "no rating change" placeholder for the two new date
# columns on every data row.
For ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# BidaskClub (row 22) has a brand new rating action on 6/16/2018 -
# record it in the newest date column and highlight it.
$ws.Range("B22").Interior.ColorIndex = 45
$ws.Range("B22").Interior.Pattern = -4142
$ws.Range("B22").Interior.Pattern = 1
$ws.Range("B22").Value = "6/16/2018,Downgrades,Hold -> Sell,"
